$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.336.78"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.742.97"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.94"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.12"
$ws.Range("E6").Value = "  -0.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.740.76"
$ws.Range("E7").Value = "  -0.52%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.43"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("E13").Value = "  -3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.09"
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.371.58"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.747.81"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.307.49"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.83"
$ws.Range("E18").Value = "  -3.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.97"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.62"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "464.32"
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.696"
$ws.Range("E23").Value = "  -1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.92"
$ws.Range("E24").Value = "  +0.88%  "

$ws.Range("E25").Value = "  +3.54%  "

$ws.Range("E26").Value = "  -1.32%  "

$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -2.46%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.891.36"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.28"
$ws.Range("E32").Value = "  -2.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.80"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.15"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.698.88"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("E38").Value = "  -1.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -5.08%  "

$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.993"
$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.01"
$ws.Range("E46").Value = "  +9.74%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.70"
$ws.Range("E47").Value = "  +3.86%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.44"
$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "388.88"
$ws.Range("E50").Value = "  -2.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.18"
$ws.Range("E51").Value = "  -0.02%  "

